# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows above the current row 18, pushing the
# existing rows 18-45 down to 21-48, then populate the 3 new rows with the
# latest "Sin especificar / Provincia de Linares" Espárragos price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:45 down by three rows (Excel default Insert copies the
# formatting -incl. the date number-format style- from the row above).
$ws.Rows(18).Resize(3).Insert()

# New row 18
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44469
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = "Espárragos"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Banquete"
$ws.Range("J18").Value = 210
$ws.Range("K18").Value = 1700
$ws.Range("L18").Value = 1800
$ws.Range("M18").Value = 1750
$ws.Range("N18").Value = "$/kilo"
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 1750
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"

# New row 19
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44469
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 300000000
$ws.Range("G19").Value = "Espárragos"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 1600
$ws.Range("M19").Value = 1550
$ws.Range("N19").Value = "$/kilo"
$ws.Range("O19").Value = "Provincia de Linares"
$ws.Range("P19").Value = 1550
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = "Hortaliza"

# New row 20
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44469
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 106
$ws.Range("K20").Value = 1400
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = 1400
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Provincia de Linares"
$ws.Range("P20").Value = 1400
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
